$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 301.25
$ws.Range("I4").Value = 211.5
$ws.Range("J4").Value = 750
$ws.Range("K4").Value = 211.5
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = -97.5
$ws.Range("N4").Value = -978
$ws.Range("H17").Value = 1251.35
$ws.Range("J17").Value = 1251.35
$ws.Range("L17").Value = 3754.05
$ws.Range("N17").Value = -4090.05
$ws.Range("H19").Value = 249.2
$ws.Range("I19").Value = 211.5
$ws.Range("K19").Value = 211.5
$ws.Range("M19").Value = -36.5
$ws.Range("H39").Value = 239.84616
$ws.Range("I39").Value = 181.14285
$ws.Range("J39").Value = 308.33334
$ws.Range("K39").Value = 543.4285500000001
$ws.Range("L39").Value = 925.0000200000001
$ws.Range("M39").Value = -247.4285500000001
$ws.Range("N39").Value = -1517.00002
$ws.Range("H40").Value = 1505.8823
$ws.Range("I40").Value = 1405
$ws.Range("J40").Value = 1650
$ws.Range("K40").Value = 1405
$ws.Range("L40").Value = 1650
$ws.Range("M40").Value = -1230
$ws.Range("N40").Value = -2000
$ws.Range("H99").Value = 705.5
$ws.Range("J99").Value = 986.6
$ws.Range("L99").Value = 2959.8
$ws.Range("N99").Value = -5955.8
$ws.Range("H100").Value = 2615
$ws.Range("I100").Value = 1376.25
$ws.Range("J100").Value = 3110.5
$ws.Range("K100").Value = 1376.25
$ws.Range("L100").Value = 3110.5
$ws.Range("M100").Value = -835.25
$ws.Range("N100").Value = -4192.5
$ws.Range("H129").Value = 981.9107
$ws.Range("J129").Value = 1093.279
$ws.Range("L129").Value = 3279.837
$ws.Range("N129").Value = -13279.837
$ws.Range("H138").Value = 2328318
$ws.Range("J138").Value = 3693.625
$ws.Range("L138").Value = 11080.875
$ws.Range("N138").Value = -21360.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1678.8
$ws.Range("I102").Value = 1598.5
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1598.5
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 23.5
$ws.Range("N102").Value = -5244
$ws.Range("H109").Value = 31377
$ws.Range("J109").Value = 31377
$ws.Range("L109").Value = 31377
$ws.Range("N109").Value = -34151

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 43250
$ws.Range("J95").Value = 43250
$ws.Range("L95").Value = 43250
$ws.Range("N95").Value = -48742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 260010.5
$ws.Range("J26").Value = 260010.5
$ws.Range("L26").Value = 260010.5
$ws.Range("N26").Value = -260584.5
$ws.Range("H31").Value = 2047.0682
$ws.Range("I31").Value = 1463.9032
$ws.Range("J31").Value = 3437.6924
$ws.Range("K31").Value = 1463.9032
$ws.Range("L31").Value = 3437.6924
$ws.Range("M31").Value = -1168.9032
$ws.Range("N31").Value = -4027.6924
$ws.Range("H34").Value = 2047.0682
$ws.Range("I34").Value = 1463.9032
$ws.Range("J34").Value = 3437.6924
$ws.Range("K34").Value = 1463.9032
$ws.Range("L34").Value = 3437.6924
$ws.Range("M34").Value = -1261.9032
$ws.Range("N34").Value = -3841.6924
$ws.Range("H141").Value = 33595
$ws.Range("J141").Value = 32314
$ws.Range("L141").Value = 32314
$ws.Range("N141").Value = -42674

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 9000
$ws.Range("N19").Value = -9348
$ws.Range("H20").Value = 1740
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 1925
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 5775
$ws.Range("M20").Value = -2773
$ws.Range("N20").Value = -6229
$ws.Range("H63").Value = 236547.11
$ws.Range("I63").Value = 302989.16
$ws.Range("K63").Value = 908967.48
$ws.Range("M63").Value = -908218.48
$ws.Range("H64").Value = 5456.303
$ws.Range("J64").Value = 6089.478
$ws.Range("L64").Value = 18268.434
$ws.Range("N64").Value = -18808.434
$ws.Range("H66").Value = 236547.11
$ws.Range("I66").Value = 302989.16
$ws.Range("K66").Value = 2726902.44
$ws.Range("M66").Value = -2723158.44
$ws.Range("H67").Value = 5456.303
$ws.Range("J67").Value = 6089.478
$ws.Range("L67").Value = 18268.434
$ws.Range("N67").Value = -20140.434
$ws.Range("H121").Value = 1293.9375
$ws.Range("J121").Value = 1368.5834
$ws.Range("L121").Value = 4105.7502
$ws.Range("N121").Value = -6725.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 27500
$ws.Range("J98").Value = 27500
$ws.Range("L98").Value = 27500
$ws.Range("N98").Value = -33490
$ws.Range("H109").Value = 30279.334
$ws.Range("J109").Value = 30279.334
$ws.Range("L109").Value = 30279.334
$ws.Range("N109").Value = -32359.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1747.8846
$ws.Range("I82").Value = 1674.8235
$ws.Range("J82").Value = 1885.8889
$ws.Range("K82").Value = 1674.8235
$ws.Range("L82").Value = 1885.8889
$ws.Range("M82").Value = -1313.8235
$ws.Range("N82").Value = -2607.8889
$ws.Range("H85").Value = 1747.8846
$ws.Range("I85").Value = 1674.8235
$ws.Range("J85").Value = 1885.8889
$ws.Range("K85").Value = 1674.8235
$ws.Range("L85").Value = 1885.8889
$ws.Range("M85").Value = -426.8235
$ws.Range("N85").Value = -4381.8889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683
$ws.Range("H95").Value = 124468.8
$ws.Range("J95").Value = 124468.8
$ws.Range("L95").Value = 124468.8
$ws.Range("N95").Value = -129960.8
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 1631.0968
$ws.Range("I122").Value = 1552.0416
$ws.Range("J122").Value = 1902.1428
$ws.Range("K122").Value = 4656.1248
$ws.Range("L122").Value = 5706.428400000001
$ws.Range("N122").Value = -10606.4284
